$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 961 already carries the per-column cell styles the new rows need (date
# format on B, centered style on D, etc.) -- copy it into each new row first so
# appended data inherits the same formatting as the rest of the table, then
# overwrite the copied values with the real data.
$templateRange = "A961:V961"
$ws.Range($templateRange).Copy($ws.Range("A962:V962"))
$ws.Range($templateRange).Copy($ws.Range("A963:V963"))
$ws.Range($templateRange).Copy($ws.Range("A964:V964"))
$ws.Range($templateRange).Copy($ws.Range("A965:V965"))
$ws.Range($templateRange).Copy($ws.Range("A966:V966"))
$ws.Range($templateRange).Copy($ws.Range("A967:V967"))
$ws.Range($templateRange).Copy($ws.Range("A968:V968"))
$ws.Range($templateRange).Copy($ws.Range("A969:V969"))
$ws.Range($templateRange).Copy($ws.Range("A970:V970"))
$ws.Range($templateRange).Copy($ws.Range("A971:V971"))
$ws.Range($templateRange).Copy($ws.Range("A972:V972"))
$ws.Range($templateRange).Copy($ws.Range("A973:V973"))
$ws.Range($templateRange).Copy($ws.Range("A974:V974"))

# Fill every column except "Type" (A) first, then "Type" last, so that brand-new
# shared strings are interned in the same order as the source workbook (the new
# match-name string is only introduced once all the new time strings exist).

# --- Row 962 (B..V) ---
$ws.Range("B962").Value2 = 45996
$ws.Range("C962").Value = "Global"
$ws.Range("D962").Value = "J-1"
$ws.Range("E962").Value = "Malik Boussaid"
$ws.Range("F962").Value = "right back"
$ws.Range("G962").Value = "01:15:48"
$ws.Range("H962").Value2 = 5.09
$ws.Range("I962").Value2 = 0.36
$ws.Range("J962").Value2 = 4.72
$ws.Range("K962").Value2 = 0.29
$ws.Range("L962").Value2 = 0.07
$ws.Range("M962").Value2 = 0.01
$ws.Range("N962").Value2 = 0
$ws.Range("O962").Value2 = 2
$ws.Range("P962").Value2 = 3.53
$ws.Range("Q962").Value2 = 25.89
$ws.Range("R962").Value2 = 4.32
$ws.Range("S962").Value2 = 39
$ws.Range("T962").Value2 = 4
$ws.Range("U962").Value2 = 29
$ws.Range("V962").Value2 = 10

# --- Row 963 (B..V) ---
$ws.Range("B963").Value2 = 45996
$ws.Range("C963").Value = "Global"
$ws.Range("D963").Value = "J-1"
$ws.Range("E963").Value = "Ilan Ihaddadene"
$ws.Range("F963").Value = "center midfield"
$ws.Range("G963").Value = "01:18:24"
$ws.Range("H963").Value2 = 5.27
$ws.Range("I963").Value2 = 0.23
$ws.Range("J963").Value2 = 5.03
$ws.Range("K963").Value2 = 0.2
$ws.Range("L963").Value2 = 0.03
$ws.Range("M963").Value2 = 0
$ws.Range("N963").Value2 = 0
$ws.Range("O963").Value2 = 0
$ws.Range("P963").Value2 = 3.94
$ws.Range("Q963").Value2 = 23.2
$ws.Range("R963").Value2 = 5.24
$ws.Range("S963").Value2 = 17
$ws.Range("T963").Value2 = 9
$ws.Range("U963").Value2 = 17
$ws.Range("V963").Value2 = 1

# --- Row 964 (B..V) ---
$ws.Range("B964").Value2 = 45997
$ws.Range("C964").Value = "Global"
$ws.Range("D964").Value = "M"
$ws.Range("E964").Value = "Emmanuel Valey"
$ws.Range("F964").Value = "left forward"
$ws.Range("G964").Value = "00:30:42"
$ws.Range("H964").Value2 = 3.96
$ws.Range("I964").Value2 = 1.04
$ws.Range("J964").Value2 = 2.9
$ws.Range("K964").Value2 = 0.59
$ws.Range("L964").Value2 = 0.3
$ws.Range("M964").Value2 = 0.15
$ws.Range("N964").Value2 = 0.02
$ws.Range("O964").Value2 = 6
$ws.Range("P964").Value2 = 7.7
$ws.Range("Q964").Value2 = 30.53
$ws.Range("R964").Value2 = 4.85
$ws.Range("S964").Value2 = 23
$ws.Range("T964").Value2 = 3
$ws.Range("U964").Value2 = 23
$ws.Range("V964").Value2 = 5

# --- Row 965 (B..V) ---
$ws.Range("B965").Value2 = 45997
$ws.Range("C965").Value = "Global"
$ws.Range("D965").Value = "M"
$ws.Range("E965").Value = "Ilan Ihaddadene"
$ws.Range("F965").Value = "center midfield"
$ws.Range("G965").Value = "01:41:20"
$ws.Range("H965").Value2 = 12.25
$ws.Range("I965").Value2 = 2.18
$ws.Range("J965").Value2 = 10.04
$ws.Range("K965").Value2 = 1.73
$ws.Range("L965").Value2 = 0.4
$ws.Range("M965").Value2 = 0.07
$ws.Range("N965").Value2 = 0
$ws.Range("O965").Value2 = 5
$ws.Range("P965").Value2 = 6.82
$ws.Range("Q965").Value2 = 29.91
$ws.Range("R965").Value2 = 5.35
$ws.Range("S965").Value2 = 44
$ws.Range("T965").Value2 = 11
$ws.Range("U965").Value2 = 26
$ws.Range("V965").Value2 = 6

# --- Row 966 (B..V) ---
$ws.Range("B966").Value2 = 45997
$ws.Range("C966").Value = "Global"
$ws.Range("D966").Value = "M"
$ws.Range("E966").Value = "Jeremie Laurent"
$ws.Range("F966").Value = "left forward"
$ws.Range("G966").Value = "00:26:28"
$ws.Range("H966").Value2 = 2.5
$ws.Range("I966").Value2 = 0.57
$ws.Range("J966").Value2 = 1.92
$ws.Range("K966").Value2 = 0.39
$ws.Range("L966").Value2 = 0.12
$ws.Range("M966").Value2 = 0.07
$ws.Range("N966").Value2 = 0
$ws.Range("O966").Value2 = 4
$ws.Range("P966").Value2 = 5.58
$ws.Range("Q966").Value2 = 29.12
$ws.Range("R966").Value2 = 4.94
$ws.Range("S966").Value2 = 12
$ws.Range("T966").Value2 = 4
$ws.Range("U966").Value2 = 14
$ws.Range("V966").Value2 = 6

# --- Row 967 (B..V) ---
$ws.Range("B967").Value2 = 45997
$ws.Range("C967").Value = "Global"
$ws.Range("D967").Value = "M"
$ws.Range("E967").Value = "Naim Dhib"
$ws.Range("F967").Value = "center midfield"
$ws.Range("G967").Value = "01:45:26"
$ws.Range("H967").Value2 = 10.77
$ws.Range("I967").Value2 = 1.81
$ws.Range("J967").Value2 = 8.94
$ws.Range("K967").Value2 = 1.47
$ws.Range("L967").Value2 = 0.32
$ws.Range("M967").Value2 = 0.05
$ws.Range("N967").Value2 = 0
$ws.Range("O967").Value2 = 5
$ws.Range("P967").Value2 = 5.99
$ws.Range("Q967").Value2 = 28.58
$ws.Range("R967").Value2 = 4.14
$ws.Range("S967").Value2 = 56
$ws.Range("T967").Value2 = 5
$ws.Range("U967").Value2 = 42
$ws.Range("V967").Value2 = 15

# --- Row 968 (B..V) ---
$ws.Range("B968").Value2 = 45997
$ws.Range("C968").Value = "Global"
$ws.Range("D968").Value = "M"
$ws.Range("E968").Value = "Malik Boussaid"
$ws.Range("F968").Value = "right back"
$ws.Range("G968").Value = "01:45:49"
$ws.Range("H968").Value2 = 12.05
$ws.Range("I968").Value2 = 2.52
$ws.Range("J968").Value2 = 9.51
$ws.Range("K968").Value2 = 1.71
$ws.Range("L968").Value2 = 0.68
$ws.Range("M968").Value2 = 0.15
$ws.Range("N968").Value2 = 0
$ws.Range("O968").Value2 = 13
$ws.Range("P968").Value2 = 6.57
$ws.Range("Q968").Value2 = 30.09
$ws.Range("R968").Value2 = 4.29
$ws.Range("S968").Value2 = 39
$ws.Range("T968").Value2 = 4
$ws.Range("U968").Value2 = 41
$ws.Range("V968").Value2 = 12

# --- Row 969 (B..V) ---
$ws.Range("B969").Value2 = 45997
$ws.Range("C969").Value = "Global"
$ws.Range("D969").Value = "M"
$ws.Range("E969").Value = "Yoan Zouma"
$ws.Range("F969").Value = "center back"
$ws.Range("G969").Value = "01:44:33"
$ws.Range("H969").Value2 = 9.14
$ws.Range("I969").Value2 = 1.15
$ws.Range("J969").Value2 = 7.97
$ws.Range("K969").Value2 = 0.76
$ws.Range("L969").Value2 = 0.32
$ws.Range("M969").Value2 = 0.07
$ws.Range("N969").Value2 = 0
$ws.Range("O969").Value2 = 5
$ws.Range("P969").Value2 = 5.07
$ws.Range("Q969").Value2 = 30.17
$ws.Range("R969").Value2 = 4.57
$ws.Range("S969").Value2 = 27
$ws.Range("T969").Value2 = 2
$ws.Range("U969").Value2 = 19
$ws.Range("V969").Value2 = 6

# --- Row 970 (B..V) ---
$ws.Range("B970").Value2 = 45997
$ws.Range("C970").Value = "Global"
$ws.Range("D970").Value = "M"
$ws.Range("E970").Value = "Amir Etien"
$ws.Range("F970").Value = "right forward"
$ws.Range("G970").Value = "01:41:44"
$ws.Range("H970").Value2 = 9.64
$ws.Range("I970").Value2 = 1.58
$ws.Range("J970").Value2 = 8.04
$ws.Range("K970").Value2 = 0.93
$ws.Range("L970").Value2 = 0.49
$ws.Range("M970").Value2 = 0.16
$ws.Range("N970").Value2 = 0.02
$ws.Range("O970").Value2 = 17
$ws.Range("P970").Value2 = 5.6
$ws.Range("Q970").Value2 = 34
$ws.Range("R970").Value2 = 5.14
$ws.Range("S970").Value2 = 37
$ws.Range("T970").Value2 = 14
$ws.Range("U970").Value2 = 33
$ws.Range("V970").Value2 = 14

# --- Row 971 (B..V) ---
$ws.Range("B971").Value2 = 45997
$ws.Range("C971").Value = "Global"
$ws.Range("D971").Value = "M"
$ws.Range("E971").Value = "Yoann Martelat"
$ws.Range("F971").Value = "center midfield"
$ws.Range("G971").Value = "01:43:39"
$ws.Range("H971").Value2 = 11.86
$ws.Range("I971").Value2 = 2.52
$ws.Range("J971").Value2 = 9.32
$ws.Range("K971").Value2 = 1.9
$ws.Range("L971").Value2 = 0.56
$ws.Range("M971").Value2 = 0.08
$ws.Range("N971").Value2 = 0
$ws.Range("O971").Value2 = 7
$ws.Range("P971").Value2 = 6.82
$ws.Range("Q971").Value2 = 28.27
$ws.Range("R971").Value2 = 4.59
$ws.Range("S971").Value2 = 39
$ws.Range("T971").Value2 = 1
$ws.Range("U971").Value2 = 37
$ws.Range("V971").Value2 = 9

# --- Row 972 (B..V) ---
$ws.Range("B972").Value2 = 45997
$ws.Range("C972").Value = "Global"
$ws.Range("D972").Value = "M"
$ws.Range("E972").Value = "Kamal Bafounta"
$ws.Range("F972").Value = "center midfield"
$ws.Range("G972").Value = "00:04:52"
$ws.Range("H972").Value2 = 0.41
$ws.Range("I972").Value2 = 0.04
$ws.Range("J972").Value2 = 0.37
$ws.Range("K972").Value2 = 0.04
$ws.Range("L972").Value2 = 0
$ws.Range("M972").Value2 = 0
$ws.Range("N972").Value2 = 0
$ws.Range("O972").Value2 = 0
$ws.Range("P972").Value2 = 3.97
$ws.Range("Q972").Value2 = 20.13
$ws.Range("R972").Value2 = 3.55
$ws.Range("S972").Value2 = 4
$ws.Range("T972").Value2 = 0
$ws.Range("U972").Value2 = 1
$ws.Range("V972").Value2 = 0

# --- Row 973 (B..V) ---
$ws.Range("B973").Value2 = 45997
$ws.Range("C973").Value = "Global"
$ws.Range("D973").Value = "M"
$ws.Range("E973").Value = "Naim Ighbane"
$ws.Range("F973").Value = "center back"
$ws.Range("G973").Value = "01:43:38"
$ws.Range("H973").Value2 = 10.31
$ws.Range("I973").Value2 = 1.38
$ws.Range("J973").Value2 = 8.91
$ws.Range("K973").Value2 = 0.93
$ws.Range("L973").Value2 = 0.38
$ws.Range("M973").Value2 = 0.09
$ws.Range("N973").Value2 = 0
$ws.Range("O973").Value2 = 7
$ws.Range("P973").Value2 = 5.92
$ws.Range("Q973").Value2 = 29.68
$ws.Range("R973").Value2 = 4.31
$ws.Range("S973").Value2 = 40
$ws.Range("T973").Value2 = 3
$ws.Range("U973").Value2 = 34
$ws.Range("V973").Value2 = 5

# --- Row 974 (B..V) ---
$ws.Range("B974").Value2 = 45997
$ws.Range("C974").Value = "Global"
$ws.Range("D974").Value = "M"
$ws.Range("E974").Value = "Sofiane Belle"
$ws.Range("F974").Value = "left forward"
$ws.Range("G974").Value = "01:19:47"
$ws.Range("H974").Value2 = 8.15
$ws.Range("I974").Value2 = 1.2
$ws.Range("J974").Value2 = 6.93
$ws.Range("K974").Value2 = 0.9
$ws.Range("L974").Value2 = 0.25
$ws.Range("M974").Value2 = 0.07
$ws.Range("N974").Value2 = 0
$ws.Range("O974").Value2 = 6
$ws.Range("P974").Value2 = 6.02
$ws.Range("Q974").Value2 = 27.95
$ws.Range("R974").Value2 = 4.53
$ws.Range("S974").Value2 = 28
$ws.Range("T974").Value2 = 4
$ws.Range("U974").Value2 = 27
$ws.Range("V974").Value2 = 4

# --- Column A ("Type") for every new row ---
$ws.Range("A962").Value = "Entrainement"
$ws.Range("A963").Value = "Entrainement"
$ws.Range("A964").Value = "N3 J10 VS Carnoux"
$ws.Range("A965").Value = "N3 J10 VS Carnoux"
$ws.Range("A966").Value = "N3 J10 VS Carnoux"
$ws.Range("A967").Value = "N3 J10 VS Carnoux"
$ws.Range("A968").Value = "N3 J10 VS Carnoux"
$ws.Range("A969").Value = "N3 J10 VS Carnoux"
$ws.Range("A970").Value = "N3 J10 VS Carnoux"
$ws.Range("A971").Value = "N3 J10 VS Carnoux"
$ws.Range("A972").Value = "N3 J10 VS Carnoux"
$ws.Range("A973").Value = "N3 J10 VS Carnoux"
$ws.Range("A974").Value = "N3 J10 VS Carnoux"

# Update the active selection to mirror where the user left off entering data
$ws.Range("E978").Select()